$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D - copy formatting from an existing header cell then set the value
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "http://dbpedia.org/ontology/party"

# Row 2 - Barack Obama
$ws.Range("A2").Value = "http://dbpedia.org/resource/Barack_Obama"
$ws.Range("B2").Value = "http://dbpedia.org/resource/1961"
$ws.Range("C2").Value = "http://dbpedia.org/resource/Honolulu"
$ws.Range("D2").Value = "http://dbpedia.org/resource/Democratic_Party_(United_States)"

# Row 3 - Benjamin Harrison
$ws.Range("A3").Value = "http://dbpedia.org/resource/Benjamin_Harrison"
$ws.Range("B3").Value = "http://dbpedia.org/resource/1833"
$ws.Range("C3").Value = "http://dbpedia.org/resource/North_Bend,_Ohio"
$ws.Range("D3").Value = "http://dbpedia.org/resource/Republican_Party_(United_States)"

# Row 4 - Calvin Coolidge
$ws.Range("A4").Value = "http://dbpedia.org/resource/Calvin_Coolidge"
$ws.Range("B4").Value = "http://dbpedia.org/resource/1872"
$ws.Range("C4").Value = "http://dbpedia.org/resource/Plymouth_Notch,_Vermont"
$ws.Range("D4").Value = "http://dbpedia.org/resource/Republican_Party_(United_States)"

# Row 5 - Harry S. Truman
$ws.Range("A5").Value = "http://dbpedia.org/resource/Harry_S._Truman"
$ws.Range("B5").Value = "http://dbpedia.org/resource/1884"
$ws.Range("C5").Value = "http://dbpedia.org/resource/Lamar"
$ws.Range("D5").Value = "http://dbpedia.org/resource/Missouri"

# Row 6 - Herbert Hoover (C6 empty)
$ws.Range("A6").Value = "http://dbpedia.org/resource/Herbert_Hoover"
$ws.Range("B6").Value = "http://dbpedia.org/resource/1874"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "http://dbpedia.org/resource/Iowa"

# Row 7 - Lyndon B. Johnson
$ws.Range("A7").Value = "http://dbpedia.org/resource/Lyndon_B._Johnson"
$ws.Range("B7").Value = "http://dbpedia.org/resource/1908"
$ws.Range("C7").Value = "http://dbpedia.org/resource/Stonewall"
$ws.Range("D7").Value = "http://dbpedia.org/resource/Texas"

# Row 8 no longer exists - clear it entirely
$ws.Range("A8:C8").Clear()
